$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.894.77"
$ws.Range("E2").Value = "  +2.61%  "

$ws.Range("D3").Value = "3.199.81"
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.96"
$ws.Range("E5").Value = "  +4.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "621.17"
$ws.Range("E6").Value = "  +1.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.386"
$ws.Range("E7").Value = "  +0.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.698"
$ws.Range("E8").Value = "  +2.58%  "

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("D10").Value = "3.197.36"
$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.570"
$ws.Range("E11").Value = "  +4.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.178"
$ws.Range("E12").Value = "  -1.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.38"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").Value = "3.789.58"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").Value = "89.686.84"
$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.66"
$ws.Range("E17").Value = "  -0.30%  "

$ws.Range("D18").Value = "3.220.24"
$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000234"
$ws.Range("E19").Value = "  +75.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.33"
$ws.Range("E20").Value = "  +11.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.31"
$ws.Range("E21").Value = "  -1.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.86"
$ws.Range("E22").Value = "  +2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.52"
$ws.Range("E23").Value = "  -0.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.03"
$ws.Range("E24").Value = "  -2.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.08"
$ws.Range("E25").Value = "  -2.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.55"
$ws.Range("E26").Value = "  -1.15%  "

$ws.Range("D27").Value = "3.356.90"
$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "75.09"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.154"
$ws.Range("E31").Value = "  -11.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.01"
$ws.Range("E32").Value = "  +33.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.40"
$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "533.95"
$ws.Range("E34").Value = "  -2.68%  "

$ws.Range("E35").Value = "  -0.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.83"
$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("E37").Value = "  -1.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.21"
$ws.Range("E38").Value = "  -0.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.30"
$ws.Range("E39").Value = "  +2.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.127"
$ws.Range("E41").Value = "  -6.06%  "

$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.370"
$ws.Range("E44").Value = "  -4.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.74"
$ws.Range("E45").Value = "  +2.44%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "171.15"
$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.20"
$ws.Range("E47").Value = "  -1.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.123"
$ws.Range("E48").Value = "  -5.19%  "

$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.23"
$ws.Range("E49").Value = "  -4.41%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.730"
$ws.Range("E50").Value = "  +3.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.612"
$ws.Range("E51").Value = "  +0.81%  "
